# Updated following manuscript revision
# The LD (light/dark) segment mapping table on Sheet1 is trimmed down to the
# first 10 "day" groups (rows 1-29) and the ld_day grouping counter in column D
# is bumped by one for the dark_pm/dark_am rows so it lines up with the
# "light" row of the same day, then the now-unused trailing rows (originally
# rows 30-58, covering days 10-20) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump column D ("ld_day") by 1 for the dark_pm / dark_am rows of each day
# so both dark rows match the following day's light row grouping id.
$rowsToBump = @(3, 4, 6, 7, 9, 10, 12, 13, 15, 16, 18, 19, 21, 22, 24, 25, 27, 28)
foreach ($r in $rowsToBump) {
    $cell = $ws.Cells.Item($r, 4)
    $cur = $cell.Value2
    $cell.Value = $cur + 1
}

# Remove the now-obsolete trailing rows (days 10-20), leaving A1:D29.
[void]$ws.Rows("30:58").Delete()

# Restore the last-saved selection/active cell shown in the sheet view.
[void]$ws.Range("G24").Select()
